$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 46072
$ws.Range("E8").Value = 147.84
$ws.Range("F8").Value = 157.84
$ws.Range("G8").Value = 147.72999999999999
$ws.Range("A9").Value = 46072
$ws.Range("E9").Value = 147.84
$ws.Range("F9").Value = 157.84
$ws.Range("G9").Value = 147.72999999999999
$ws.Range("A10").Value = 46072
$ws.Range("D10").Value = 159.18
$ws.Range("E10").Value = 150.56
$ws.Range("F10").Value = 160.56
$ws.Range("G10").Value = 150.80000000000001
$ws.Range("A11").Value = 46071
$ws.Range("D11").Value = 157.63
$ws.Range("E11").Value = 148.41999999999999
$ws.Range("F11").Value = 158.41999999999999
$ws.Range("G11").Value = 148.31
$ws.Range("A12").Value = 46071
$ws.Range("D12").Value = 157.63
$ws.Range("E12").Value = 148.41999999999999
$ws.Range("F12").Value = 158.41999999999999
$ws.Range("G12").Value = 148.31
$ws.Range("A13").Value = 46071
$ws.Range("D13").Value = 159.27000000000001
$ws.Range("E13").Value = 151.13
$ws.Range("F13").Value = 161.13
$ws.Range("G13").Value = 151.37
$ws.Range("A17").Value = 46072
$ws.Range("D17").Value = 163.58000000000001
$ws.Range("E17").Value = 154.13
$ws.Range("F17").Value = 164.13
$ws.Range("A18").Value = 46071
$ws.Range("D18").Value = 163.65
$ws.Range("E18").Value = 154.69
$ws.Range("F18").Value = 164.69
$ws.Range("A22").Value = 46072
$ws.Range("D22").Value = 158.91999999999999
$ws.Range("E22").Value = 150.18
$ws.Range("F22").Value = 159.78
$ws.Range("G22").Value = 151.93
$ws.Range("A23").Value = 46072
$ws.Range("D23").Value = 164.15
$ws.Range("E23").Value = 156.26
$ws.Range("F23").Value = 166.26
$ws.Range("A24").Value = 46072
$ws.Range("D24").Value = 164.34
$ws.Range("E24").Value = 156.79
$ws.Range("F24").Value = 166.79
$ws.Range("A25").Value = 46072
$ws.Range("D25").Value = 164.35
$ws.Range("E25").Value = 156.30000000000001
$ws.Range("F25").Value = 166.3
$ws.Range("G25").Value = 157.15
$ws.Range("A26").Value = 46072
$ws.Range("D26").Value = 163.98
$ws.Range("E26").Value = 157.88
$ws.Range("F26").Value = 167.88
$ws.Range("A27").Value = 46071
$ws.Range("D27").Value = 159.03
$ws.Range("E27").Value = 150.76
$ws.Range("F27").Value = 160.36000000000001
$ws.Range("G27").Value = 152.52000000000001
$ws.Range("A28").Value = 46071
$ws.Range("D28").Value = 164.24
$ws.Range("E28").Value = 156.83000000000001
$ws.Range("F28").Value = 166.83
$ws.Range("A29").Value = 46071
$ws.Range("D29").Value = 164.43
$ws.Range("E29").Value = 157.36000000000001
$ws.Range("F29").Value = 167.36
$ws.Range("A30").Value = 46071
$ws.Range("D30").Value = 164.44
$ws.Range("E30").Value = 156.87
$ws.Range("F30").Value = 166.87
$ws.Range("G30").Value = 157.72999999999999
$ws.Range("A31").Value = 46071
$ws.Range("D31").Value = 164.07
$ws.Range("E31").Value = 158.46
$ws.Range("F31").Value = 168.46
$ws.Range("A35").Value = 46072
$ws.Range("D35").Value = 157.69
$ws.Range("E35").Value = 148.26
$ws.Range("F35").Value = 157.26
$ws.Range("A36").Value = 46071
$ws.Range("D36").Value = 157.77000000000001
$ws.Range("E36").Value = 148.82
$ws.Range("F36").Value = 157.82
$ws.Range("A40").Value = 46072
$ws.Range("D40").Value = 163.86
$ws.Range("E40").Value = 155.37
$ws.Range("F40").Value = 165.37
$ws.Range("A41").Value = 46072
$ws.Range("D41").Value = 163.58000000000001
$ws.Range("E41").Value = 155.79
$ws.Range("F41").Value = 165.79
$ws.Range("A42").Value = 46071
$ws.Range("D42").Value = 163.96
$ws.Range("E42").Value = 156.03
$ws.Range("F42").Value = 166.03
$ws.Range("A43").Value = 46071
$ws.Range("D43").Value = 163.68
$ws.Range("E43").Value = 156.46
$ws.Range("F43").Value = 166.46
$ws.Range("A47").Value = 46072
$ws.Range("D47").Value = 158.6
$ws.Range("E47").Value = 150.54
$ws.Range("F47").Value = 160.54
$ws.Range("A48").Value = 46072
$ws.Range("D48").Value = 158.30000000000001
$ws.Range("E48").Value = 150.53
$ws.Range("F48").Value = 160.53
$ws.Range("A49").Value = 46071
$ws.Range("D49").Value = 159.01
$ws.Range("E49").Value = 150.65
$ws.Range("F49").Value = 160.65
$ws.Range("A50").Value = 46071
$ws.Range("D50").Value = 158.71
$ws.Range("E50").Value = 150.63
$ws.Range("F50").Value = 160.63
$ws.Range("A54").Value = 46072
$ws.Range("D54").Value = 173.06
$ws.Range("E54").Value = 163.08000000000001
$ws.Range("F54").Value = 173.08
$ws.Range("A55").Value = 46072
$ws.Range("D55").Value = 162.43
$ws.Range("E55").Value = 162.28
$ws.Range("F55").Value = 172.28
$ws.Range("A56").Value = 46072
$ws.Range("D56").Value = 162.19
$ws.Range("A57").Value = 46072
$ws.Range("D57").Value = 163.11000000000001
$ws.Range("E57").Value = 156.69999999999999
$ws.Range("A58").Value = 46072
$ws.Range("D58").Value = 158.88
$ws.Range("E58").Value = 152.6
$ws.Range("F58").Value = 162.6
$ws.Range("A59").Value = 46072
$ws.Range("D59").Value = 166.02
$ws.Range("E59").Value = 161.65
$ws.Range("A60").Value = 46071
$ws.Range("D60").Value = 173.14
$ws.Range("E60").Value = 163.69999999999999
$ws.Range("F60").Value = 173.7
$ws.Range("A61").Value = 46071
$ws.Range("D61").Value = 162.52000000000001
$ws.Range("E61").Value = 162.72
$ws.Range("F61").Value = 172.72
$ws.Range("A62").Value = 46071
$ws.Range("D62").Value = 162.28
$ws.Range("A63").Value = 46071
$ws.Range("D63").Value = 163.18
$ws.Range("E63").Value = 157.13999999999999
$ws.Range("A64").Value = 46071
$ws.Range("D64").Value = 158.94999999999999
$ws.Range("E64").Value = 153.04
$ws.Range("F64").Value = 163.04
$ws.Range("A65").Value = 46071
$ws.Range("D65").Value = 166.09
$ws.Range("E65").Value = 162.25
